$wb = $excel.ActiveWorkbook

# Sheet names (in tab order) and the edits each one needs:
#   1: Potencia Acumulada - SIN (MW)    -> E1 label "2050",      drop Total row 13
#   2: Geracao Periodo Medio (MWMed)    -> E1 label "2050",      drop Total row 13
#   3: Atendimento a Ponta(MW)          -> E1 label "2050",      drop Total row 13
#   4: Potencia Incremental - SIN(MW)   -> E1 label "2041-2050", drop Total row 13
#   5: Emissoes Totais (MtCO2eq)        -> E1 label "2050"       (no Total row here)
#   6: Custo Total (bilhões de R$)      -> drop Total row 4

function Set-HeaderLabel {
    param($ws, [string]$cellRef, [string]$sourceCellRef, [string]$text)
    # A bare Range.Value = "2050" gets auto-coerced to a NUMBER (just like
    # typing it into Excel would), which is not what we want -- the column
    # header must stay a text label, matching the other header cells.
    # Prefixing with an apostrophe forces text, but Excel then stamps the
    # cell with a "quote prefix" style variant instead of the original
    # bold/centered/bordered header style. Re-apply that original style
    # (taken from the neighboring header cell, e.g. D1) via PasteSpecial so
    # the cell ends up identical in formatting to the other header cells.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($sourceCellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$ws1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-HeaderLabel $ws1 "E1" "D1" "2050"
$ws1.Rows.Item(13).Delete()

$ws2 = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-HeaderLabel $ws2 "E1" "D1" "2050"
$ws2.Rows.Item(13).Delete()

$ws3 = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-HeaderLabel $ws3 "E1" "D1" "2050"
$ws3.Rows.Item(13).Delete()

$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-HeaderLabel $ws4 "E1" "D1" "2041-2050"
$ws4.Rows.Item(13).Delete()

$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-HeaderLabel $ws5 "E1" "D1" "2050"

$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()

$excel.CutCopyMode = $false
